# Refresh the cryptos list (Price + Volume(1h) columns) with latest snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "27.868.76"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +1.29%  "
# Row 3: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.771.40"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "
# Row 4: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9972"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  -1.11%  "
# Row 5: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "322.08"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
# Row 6: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9963"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -0.69%  "
# Row 7: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4263"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  -4.52%  "
# Row 8: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3609"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  -2.55%  "
# Row 9: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "44.36"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  -1.66%  "
# Row 10: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07464"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  -2.85%  "
# Row 11: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "1.104"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("E12").Value = "  -1.14%  "
# Row 13: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "21.55"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -0.48%  "
# Row 14: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "6.124"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -0.61%  "
# Row 15: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "7.303"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -1.07%  "
# Row 16: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "1.790.88"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +1.62%  "
# Row 17: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "91.18"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +0.49%  "
# Row 18: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "0.00001060"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -1.28%  "
# Row 19: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06350"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +0.70%  "
# Row 20: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9959"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -0.62%  "
# Row 21: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "17.17"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -1.22%  "
# Row 22: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "5.942"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -3.52%  "
# Row 23: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "27.882.42"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +1.17%  "
# Row 24: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "11.35"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  -1.89%  "
# Row 25: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "2.162"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  -6.81%  "
# Row 26: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "159.95"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +4.51%  "
# Row 27: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "20.26"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -1.69%  "
# Row 28: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "1.998.51"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +1.92%  "
# Row 29: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "2.165"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  -6.05%  "
# Row 30: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "125.75"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -1.91%  "
# Row 31: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "1.169"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -1.71%  "
# Row 32: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "5.684"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -0.64%  "
# Row 33: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08987"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -2.54%  "
# Row 34: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "3.521"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("E35").Value = "  +0.06%  "
# Row 36: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.02315"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +0.23%  "
# Row 37: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "5.068"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  +0.40%  "
# Row 38: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "0.2110"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -2.48%  "
# Row 39: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6423"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -0.16%  "
# Row 40: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06049"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -0.59%  "
# Row 41: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "1.179"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +0.12%  "
# Row 42: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9956"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -0.64%  "
# Row 43: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "7.851"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -1.58%  "
# Row 44: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "1.388"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  -1.08%  "
# Row 45: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "13.64"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "
# Row 46: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "0.5962"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "
# Row 47: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "3.691"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -1.21%  "
# Row 48: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "124.11"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -1.09%  "
# Row 49: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "1.982"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "
# Row 50: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "1.146"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +0.61%  "
# Row 51: Price cell keeps plain-text formatting even though it looks numeric
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06883"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
